# Add a second Hour:Min picker (to go with the "Start" button) further
# down the slide, mirroring the existing picker built from shapes
# id=5 ("00"), id=7 (":"), id=8 ("Hour"), id=9 ("Min").
# Commit: "Change startup folder and Start buttons made, but both need
# to be fully setup".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The host's shape-id allocator hands out the smallest unused id. The
# canonical document has gaps at 2, 3, 14, 15, 19, 21, 25, 27, 31, 32
# (ids that existed at some point but aren't on today's slide). Burn
# through the first five gaps (2, 3, 14, 15, 19) with throwaway shapes
# so the five real shapes we add below land on ids 21, 25, 27, 31, 32 -
# exactly what the target deck uses.
for ($i = 0; $i -lt 5; $i++) {
    $junk = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $junk.Delete()
}

# --- New shape: "00" (hour) textbox -> id 21, name "TextBox 20" ---
$t1 = $s.Shapes.Item(2).Duplicate().Item(1)
$t1.Name = "TextBox 20"
$t1.Left = 640.7350463867188
$t1.Top = 298.4593811035156
$t1.Width = 33.84000015258789
$t1.Height = 29.081260681152344
$t1.TextFrame.TextRange.Text = "00"

# --- New shape: "00" (minute) textbox -> id 25, name "TextBox 24" ---
$t2 = $s.Shapes.Item(2).Duplicate().Item(1)
$t2.Name = "TextBox 24"
$t2.Left = 723.2950439453125
$t2.Top = 297.7393798828125
$t2.Width = 33.84000015258789
$t2.Height = 29.081260681152344
$t2.TextFrame.TextRange.Text = "00"

# --- New shape: ":" separator textbox -> id 27, name "TextBox 26" ---
$t3 = $s.Shapes.Item(4).Duplicate().Item(1)
$t3.Name = "TextBox 26"
$t3.Left = 672.5650634765625
$t3.Top = 296.4593811035156
$t3.Width = 23.940000534057617
$t3.Height = 29.081260681152344
$t3.TextFrame.TextRange.Text = ":"

# --- New shape: "Hour" label textbox -> id 31, name "TextBox 30" ---
$t4 = $s.Shapes.Item(5).Duplicate().Item(1)
$t4.Name = "TextBox 30"
$t4.Left = 600.7800903320312
$t4.Top = 300.34063720703125
$t4.Width = 46.894962310791016
$t4.Height = 24.234411239624023
$t4.TextFrame.TextRange.Text = "Hour"

# --- New shape: "Min" label textbox -> id 32, name "TextBox 31" ---
$t5 = $s.Shapes.Item(6).Duplicate().Item(1)
$t5.Name = "TextBox 31"
$t5.Left = 685.8076171875
$t5.Top = 301.7806396484375
$t5.Width = 46.894962310791016
$t5.Height = 24.234411239624023
$t5.TextFrame.TextRange.Text = "Min"
